# Update the "TUERCAS" price list sheet:
#  - bump the date in A1 by one month (30-day price list refresh)
#  - refresh the zinc-plated nut prices in column D for the two price tables

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header date: 24/04/2024 -> 24/05/2024
$ws.Range("A1").Value = 45436

# "Tuercas altas" price table (rows 23-28)
$ws.Range("D23").Value = 11050
$ws.Range("D24").Value = 13650
$ws.Range("D25").Value = 19500
$ws.Range("D26").Value = 28600
$ws.Range("D27").Value = 52000
$ws.Range("D28").Value = 71500

# "Tuercas bajas" price table (rows 36-37)
$ws.Range("D36").Value = 5920
$ws.Range("D37").Value = 8730

# Re-assert the A30:D30 merge last so it is re-emitted after A1:E1 in the
# saved mergeCells list (matches canonical ordering after the refresh).
$ws.Range("A30:D30").UnMerge()
$ws.Range("A30:D30").Merge()
